# Updates cryptos list values (price + 1h volume change) per the
# "Updated cryptos list on Sun Oct 22 13:31:24 UTC 2023 with GitHub Actions" commit.
# Rows 22/23 (Avalanche/Uniswap) also swap places in this update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.948.94"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").Value = "1.634.14"
$ws.Range("E3").Value = "  +1.89%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'214.76"
$ws.Range("E5").Value = "  +1.21%  "

$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").Value = "'28.84"
$ws.Range("E8").Value = "  -0.34%  "

$ws.Range("E9").Value = "  +0.49%  "

$ws.Range("E10").Value = "  +0.26%  "

$ws.Range("D11").Value = "'0.0904"
$ws.Range("E11").Value = "  -0.25%  "

$ws.Range("D12").Value = "1.866.30"
$ws.Range("E12").Value = "  +1.81%  "

$ws.Range("D13").Value = "1.634.05"
$ws.Range("E13").Value = "  +2.20%  "

$ws.Range("E14").Value = "  +1.03%  "

$ws.Range("D15").Value = "'9.30"
$ws.Range("E15").Value = "  +14.41%  "

$ws.Range("D16").Value = "29.967.06"
$ws.Range("E16").Value = "  +0.88%  "

$ws.Range("D17").Value = "'3.85"
$ws.Range("E17").Value = "  +1.35%  "

$ws.Range("D18").Value = "'64.19"
$ws.Range("E18").Value = "  -0.13%  "

$ws.Range("D19").Value = "'241.10"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "0.0₃0703"
$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'4.14"
$ws.Range("E22").Value = "  +2.59%  "

$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "'9.81"
$ws.Range("E23").Value = "  +3.35%  "

$ws.Range("D24").Value = "'2.16"
$ws.Range("E24").Value = "  +2.88%  "

$ws.Range("D25").Value = "'157.74"
$ws.Range("E25").Value = "  +1.03%  "

$ws.Range("D26").Value = "'15.50"
$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("D27").Value = "'0.110"
$ws.Range("E27").Value = "  +0.57%  "

$ws.Range("D28").Value = "'6.58"
$ws.Range("E28").Value = "  +1.46%  "

$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("E30").Value = "  +2.46%  "

$ws.Range("E31").Value = "  +3.71%  "

$ws.Range("D32").Value = "'3.40"
$ws.Range("E32").Value = "  +4.67%  "

$ws.Range("E33").Value = "  +0.37%  "

$ws.Range("D34").Value = "1.431.85"
$ws.Range("E34").Value = "  +0.56%  "

$ws.Range("E35").Value = "  +5.26%  "

$ws.Range("E36").Value = "  -1.43%  "

$ws.Range("D37").Value = "'2.78"
$ws.Range("E37").Value = "  -2.79%  "

$ws.Range("E38").Value = "  -0.42%  "

$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").Value = "'76.04"
$ws.Range("E40").Value = "  +12.24%  "

$ws.Range("E41").Value = "  +0.58%  "

$ws.Range("D42").Value = "'2.00"
$ws.Range("E42").Value = "  +1.76%  "

$ws.Range("D43").Value = "'0.831"
$ws.Range("E43").Value = "  +1.48%  "

$ws.Range("E44").Value = "  -0.47%  "

$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("E46").Value = "  +0.62%  "

$ws.Range("D47").Value = "'51.41"
$ws.Range("E47").Value = "  -7.37%  "

$ws.Range("D48").Value = "'5.36"
$ws.Range("E48").Value = "  -1.29%  "

$ws.Range("D49").Value = "1.773.59"
$ws.Range("E49").Value = "  +1.98%  "

$ws.Range("E50").Value = "  +13.17%  "

$ws.Range("D51").Value = "'90.44"
$ws.Range("E51").Value = "  +4.22%  "
